# Generate Report for Handoff
# Status moves from "In Translation" to "Ready for handoff", and the
# handoff/generation timestamps are refreshed on the Overview sheet and
# the two per-locale sheets (zh-cn, de-de).

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# --- Overview sheet -------------------------------------------------
# E2/F2 hold the per-locale status ("zh-cn" / "de-de" columns)
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
# G2 holds "Latest HO Xliff Generate Date"
$wsOverview.Range("G2").Value = "2016-11-02 04:51:00"

# --- zh-cn sheet ------------------------------------------------------
# C2 holds "Status"
$wsZhCn.Range("C2").Value = "Ready for handoff"
# H2 holds "Latest Handoff Datetime"
$wsZhCn.Range("H2").Value = "2016-11-02 04:50:46"

# --- de-de sheet ------------------------------------------------------
# C2 holds "Status"
$wsDeDe.Range("C2").Value = "Ready for handoff"
# H2 holds "Latest Handoff Datetime"
$wsDeDe.Range("H2").Value = "2016-11-02 04:51:00"

# --- Column widths ------------------------------------------------------
# The Status columns widened (content grew from "In Translation" to the
# longer "Ready for handoff"), matching an autofit-style relayout.
$wsOverview.Columns.Item(5).ColumnWidth = 16.25
$wsOverview.Columns.Item(6).ColumnWidth = 16.25
$wsZhCn.Columns.Item(3).ColumnWidth = 16.25
$wsDeDe.Columns.Item(3).ColumnWidth = 16.25
